$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing columns
# A:E (name, national id, department, amount, signature) to B:F and
# carries their column widths / cell styles along automatically.
$ws.Columns("A").Insert()

# Fill the new column A with sequential row numbers for each data row
# (rows 2-34), matching the style (border / bold 18pt font / centered,
# shrink-to-fit alignment / 0.00 number format) already used by the rest
# of that row. The header row (1) and the "total" row (35) are left
# without a value in column A.
for ($i = 2; $i -le 34; $i++) {
    $cell = $ws.Cells.Item($i, 1)
    $cell.Value = $i - 1
    $cell.Font.Bold = $true
    $cell.Font.Size = 18
    $cell.NumberFormat = "0.00"
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    $cell.ShrinkToFit = $true
}
